$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "goalkeeper pants youth soccer"
$ws.Range("A2").Value = "compression pants men cheap"
$ws.Range("A3").Value = "leggings baseball"
$ws.Range("A4").Value = "indoor knee pads"
$ws.Range("A5").Value = "yoga knee pads 2 pack"
$ws.Range("A6").Value = "compression shorts for basketball"
$ws.Range("A7").Value = "knee pads gym"
$ws.Range("A8").Value = "boys basketball leggings youth"
$ws.Range("A9").Value = "knee sweat pants"
$ws.Range("A10").Value = "knee pads for sports"
$ws.Range("A11").Value = "black softball pants"
$ws.Range("A12").Value = "arthritis hope knee compression sleeve"
$ws.Range("A13").Value = "honeycomb tights"
$ws.Range("A14").Value = "mens sheer pants"
$ws.Range("A15").Value = "mens basketball knee sleeves"
$ws.Range("A16").Value = "boys baseball pants short"
$ws.Range("A17").Value = "spandex work pants men"
$ws.Range("A18").Value = "baseball pants long"
$ws.Range("A19").Value = "capri pants for men"
$ws.Range("A20").Value = "boys basketball tights"
$ws.Range("A21").Value = "tight pants"
$ws.Range("A22").Value = "knee length shorts for men"
$ws.Range("A23").Value = "softball shorts for men"
$ws.Range("A24").Value = "baseball catchers hand pad"
$ws.Range("A25").Value = "cold pad knee"
$ws.Range("A26").Value = "football pants mens with pads"
$ws.Range("A27").Value = "yoga positions chart"
$ws.Range("A28").Value = "hex leg sleeve youth"
$ws.Range("A29").Value = "knees pads yoga"
$ws.Range("A30").Value = "work pants for men with knee pad"
$ws.Range("A31").Value = "tight basketball shorts"
$ws.Range("A32").Value = "adult hockey pants"
$ws.Range("A33").Value = "youth sliding shorts baseball"
$ws.Range("A34").Value = "sheer pants for men"
$ws.Range("A35").Value = "paintball shorts"
$ws.Range("A36").Value = "lacrosse padded shorts"
$ws.Range("A37").Value = "knee pads for paintball"
$ws.Range("A38").Value = "mens big and tall basketball pants"
$ws.Range("A39").Value = "mens shorts below the knee"
$ws.Range("A40").Value = "basketball tights youth"
$ws.Range("A41").Value = "girls sliding shorts softball"
$ws.Range("A42").Value = "youth compression pants boys"
$ws.Range("A43").Value = "boys compression running pants"
$ws.Range("A44").Value = "compression yoga tights"
$ws.Range("A45").Value = "sliding shorts youth girls softball"
$ws.Range("A46").Value = "basketball lot"
$ws.Range("A47").Value = "paintball leg pads"
$ws.Range("A48").Value = "adult black football pants"
$ws.Range("A49").Value = "knee hockey pads"
$ws.Range("A50").Value = "knee sleeve basketball"
$ws.Range("A51").Value = "knee pads work pants"
$ws.Range("A52").Value = "mens knee pads"
$ws.Range("A53").Value = "padded baseball sliding shorts"
$ws.Range("A54").Value = "soccer tights for men"
$ws.Range("A55").Value = "boys tight pants"
$ws.Range("A56").Value = "black leggings for men"
$ws.Range("A57").Value = "men athletic compression pants"
$ws.Range("A58").Value = "youth 5 pad girdle"
$ws.Range("A59").Value = "poc knee pads"
$ws.Range("A60").Value = "men baseball pants black"
$ws.Range("A61").Value = "youth small compression pants"
$ws.Range("A62").Value = "hex knee pads compression leg sleeve"
$ws.Range("A63").Value = "hockey compression pants"
$ws.Range("A64").Value = "softball mens"
$ws.Range("A65").Value = "x compression pants"
$ws.Range("A66").Value = "men compression legging"
$ws.Range("A67").Value = "wrestling knee pads"
$ws.Range("A68").Value = "boys tights for sports youth"
$ws.Range("A69").Value = "hockey pants youth"
$ws.Range("A70").Value = "pant sport men"
$ws.Range("A71").Value = "girls softball sliding pants"
$ws.Range("A72").Value = "boys small compression pants"
$ws.Range("A73").Value = "compression spandex for men"
$ws.Range("A74").Value = "knee pads mtb"
$ws.Range("A75").Value = "sport leggings men"
$ws.Range("A76").Value = "mens sliding shorts"
$ws.Range("A77").Value = "kneeling on the promises"
$ws.Range("A78").Value = "tall baseball pants mens"
$ws.Range("A79").Value = "yoga knee pad thick"
$ws.Range("A80").Value = "men gym leggings"
$ws.Range("A81").Value = "basketball knee sleeves for men"
$ws.Range("A82").Value = "basketball knee sleeves with pads"
$ws.Range("A83").Value = "youth athletic tights"
$ws.Range("A84").Value = "baseball pants adults"
$ws.Range("A85").Value = "tights for soccer"
$ws.Range("A86").Value = "sliding shorts"
$ws.Range("A87").Value = "mens sport leggings"
$ws.Range("A88").Value = "6 inch basketball"
$ws.Range("A89").Value = "football waist pads"
$ws.Range("A90").Value = "mens work pants knee pads"
$ws.Range("A91").Value = "knee sleeve lacrosse"
$ws.Range("A92").Value = "men sliding shorts"
$ws.Range("A93").Value = "padded leggings"
$ws.Range("A94").Value = "yoga pads"
$ws.Range("A95").Value = "mens cold tights"
$ws.Range("A96").Value = "boy tights youth"
$ws.Range("A97").Value = "basketball cycling"
$ws.Range("A98").Value = "basketball knee sleeves youth"
$ws.Range("A99").Value = "soccer compression pants"
$ws.Range("A100").Value = "sliding pads"
